# "Generate Report for Handback" -- localization-status.xlsx
#
# For each language sheet (zh-cn, de-de) the two data rows have now been
# handed back to the source: the Status column flips from "Ready for
# handoff" to "Handed back: in sync with en-US" (this text is a shared
# string also used, unchanged, by the Overview sheet, so replacing it
# everywhere keeps Overview in sync too); the "Latest Target File" /
# "Latest Handback File" columns (F/G) get populated with hyperlinked
# filenames; and the "Latest Handback DateTime" column (H) moves from the
# zero-date placeholder to the real handback timestamp.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# Flip the status text everywhere it appears (Overview + both language
# sheets) in one go -- this keeps all cells referencing the shared string
# in sync without having to touch the Overview sheet directly.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($statusOld, $statusNew)
}

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- zh-cn sheet -----------------------------------------------------
$zhcn.Range("H2").Value2 = "2016-03-21 22:45:04"
$zhcn.Range("H3").Value2 = "2016-03-21 22:45:04"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/bb2163d0345e828b580b99bbe2e926741744e141/e2e/3ecf4bd1-77c9-4272-a80d-af88cb797606.md",
    "",
    "",
    "3ecf4bd1-77c9-4272-a80d-af88cb797606.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/852a125ed5938c670ca7030f1d206ba7750cce0b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3ecf4bd1-77c9-4272-a80d-af88cb797606.9d4d96d8a27b469758264c76051ffcb00d4002cd.zh-cn.xlf",
    "",
    "",
    "3ecf4bd1-77c9-4272-a80d-af88cb797606.9d4d96d8a27b469758264c76051ffcb00d4002cd.zh-cn.xlf"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/bb2163d0345e828b580b99bbe2e926741744e141/e2e/9b2c1858-178e-4bff-b5ee-50b671b78afc.md",
    "",
    "",
    "9b2c1858-178e-4bff-b5ee-50b671b78afc.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/852a125ed5938c670ca7030f1d206ba7750cce0b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9b2c1858-178e-4bff-b5ee-50b671b78afc.ed2a2363bfd6084c5a3da6fffc05efa45ab38861.zh-cn.xlf",
    "",
    "",
    "9b2c1858-178e-4bff-b5ee-50b671b78afc.ed2a2363bfd6084c5a3da6fffc05efa45ab38861.zh-cn.xlf"
)

# --- de-de sheet -------------------------------------------------------
$dede.Range("H2").Value2 = "2016-03-21 22:45:12"
$dede.Range("H3").Value2 = "2016-03-21 22:45:12"

$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/bb2163d0345e828b580b99bbe2e926741744e141/e2e/3ecf4bd1-77c9-4272-a80d-af88cb797606.md",
    "",
    "",
    "3ecf4bd1-77c9-4272-a80d-af88cb797606.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/51212d6e2425f8284484d65e4d8a0d5b4048be8c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3ecf4bd1-77c9-4272-a80d-af88cb797606.9d4d96d8a27b469758264c76051ffcb00d4002cd.de-de.xlf",
    "",
    "",
    "3ecf4bd1-77c9-4272-a80d-af88cb797606.9d4d96d8a27b469758264c76051ffcb00d4002cd.de-de.xlf"
)
$dede.Hyperlinks.Add(
    $dede.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/bb2163d0345e828b580b99bbe2e926741744e141/e2e/9b2c1858-178e-4bff-b5ee-50b671b78afc.md",
    "",
    "",
    "9b2c1858-178e-4bff-b5ee-50b671b78afc.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/51212d6e2425f8284484d65e4d8a0d5b4048be8c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9b2c1858-178e-4bff-b5ee-50b671b78afc.ed2a2363bfd6084c5a3da6fffc05efa45ab38861.de-de.xlf",
    "",
    "",
    "9b2c1858-178e-4bff-b5ee-50b671b78afc.ed2a2363bfd6084c5a3da6fffc05efa45ab38861.de-de.xlf"
)

Write-Output "Handback report generated"
